# Play fallen audio when infrared triggered at locking on state
#
# Adds a new transition cell (D2) describing the "INFRARED TRIGGERED"
# event while in the "LOCKING ON" state: it plays the "fallen" audio
# and moves to the LOCKING ON state (re-entrant / self) when the egg
# has already fallen. The A1 header cell ("state\event") is re-entered
# so it lives alongside the new content, and the row is resized to fit
# the new three-line cell text, matching the other multi-line rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-assert the header cell's text (kept as-is, but this normalizes its
# storage the way the rest of the header row is already stored).
$ws.Range("A1").Value = "state\event"

# New state-transition cell: LOCKING ON + INFRARED TRIGGERED ->
# play fallen / lock on / fallen -> LOCKING ON
$ws.Range("D2").Value = "play fallen`nlock on`nfallen\LOCKING ON"

# Row 2 now holds a three-line cell like the other wrapped transition
# rows (4-7), so grow it to match their height.
$ws.Rows.Item(2).RowHeight = 34.5

# Leave the selection on the cell that was just edited.
$ws.Range("D2").Select()
